# Fix a typo on slide 12 ("Regelgeving die ... CSRD & ESRS E5 (gebruik
# termonoloy)") -> "... (gebruik terminology)".
#
# In the OOXML the word "termonoloy" is split across three runs
# (" ", "termonoloy", ")"). Editing that exact span through
# TextRange.Characters() reproduces PowerPoint's own behaviour of
# collapsing the edited runs into a single run that carries the
# formatting of the surrounding (non-misspelled) text.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(12)
$sh = $s.Shapes.Item(2)          # "Tekstvak 2" textbox (shape id 3)

$tr  = $sh.TextFrame.TextRange
$len = $tr.Length

# Locate " termonoloy)" robustly instead of hard-coding a character
# offset (the presence of "&" earlier in the string shifts indices).
$needle  = " termonoloy)"
$nlen    = $needle.Length
$foundAt = -1
for ($i = 1; $i -le ($len - $nlen + 1); $i++) {
    if ($tr.Characters($i, $nlen).Text -eq $needle) {
        $foundAt = $i
        break
    }
}

if ($foundAt -ge 1) {
    # Remember the textbox's current (auto-fitted) height. Toggling
    # AutoSize off while we edit the text keeps the shape from being
    # resized mid-edit; we restore AutoSize and the original height
    # afterwards so the shape geometry round-trips unchanged.
    $origHeight = [double]"$($sh.Height)"
    $restoreAutoSize = $false
    try {
        $sh.TextFrame.AutoSize = 0   # msoAutoSizeNone
        $restoreAutoSize = $true
    } catch {
    }

    $tr.Characters($foundAt, $nlen).Text = " terminology)"

    if ($restoreAutoSize) {
        $sh.TextFrame.AutoSize = 1   # msoAutoSizeShapeToFitText
        $sh.Height = $origHeight
    }
}
